{"js": "const replacements = [\n  [\"275\u00f77=\", \"714\u00f76=\"],\n  [\"621\u00f73=\", \"137\u00f79=\"],\n  [\"655\u00f75=\", \"666\u00f72=\"],\n  [\"412\u00f78=\", \"712\u00f79=\"],\n  [\"323\u00f78=\", \"661\u00f73=\"],\n  [\"542\u00f72=\", \"675\u00f75=\"],\n  [\"732\u00f78=\", \"205\u00f74=\"],\n  [\"829\u00f75=\", \"806\u00f75=\"],\n  [\"215\u00f77=\", \"528\u00f74=\"],\n  [\"723\u00f79=\", \"258\u00f73=\"],\n  [\"793\u00f73=\", \"605\u00f75=\"],\n  [\"602\u00f75=\", \"352\u00f73=\"],\n  [\"721\u00f78=\", \"459\u00f79=\"],\n  [\"640\u00f74=\", \"880\u00f77=\"],\n  [\"674\u00f79=\", \"774\u00f74=\"],\n  [\"146\u00f79=\", \"940\u00f77=\"],\n  [\"471\u00f79=\", \"325\u00f77=\"],\n  [\"162\u00f75=\", \"962\u00f72=\"],\n  [\"201\u00f74=\", \"234\u00f79=\"],\n  [\"736\u00f73=\", \"604\u00f76=\"],\n  [\"502\u00f77=\", \"640\u00f76=\"],\n  [\"619\u00f77=\", \"845\u00f77=\"],\n  [\"628\u00f73=\", \"849\u00f74=\"],\n  [\"641\u00f72=\", \"363\u00f78=\"],\n  [\"284\u00f74=\", \"332\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('275\u00f77=', '714\u00f76='),\n    @('621\u00f73=', '137\u00f79='),\n    @('655\u00f75=', '666\u00f72='),\n    @('412\u00f78=', '712\u00f79='),\n    @('323\u00f78=', '661\u00f73='),\n    @('542\u00f72=', '675\u00f75='),\n    @('732\u00f78=', '205\u00f74='),\n    @('829\u00f75=', '806\u00f75='),\n    @('215\u00f77=', '528\u00f74='),\n    @('723\u00f79=', '258\u00f73='),\n    @('793\u00f73=', '605\u00f75='),\n    @('602\u00f75=', '352\u00f73='),\n    @('721\u00f78=', '459\u00f79='),\n    @('640\u00f74=', '880\u00f77='),\n    @('674\u00f79=', '774\u00f74='),\n    @('146\u00f79=', '940\u00f77='),\n    @('471\u00f79=', '325\u00f77='),\n    @('162\u00f75=', '962\u00f72='),\n    @('201\u00f74=', '234\u00f79='),\n    @('736\u00f73=', '604\u00f76='),\n    @('502\u00f77=', '640\u00f76='),\n    @('619\u00f77=', '845\u00f77='),\n    @('628\u00f73=', '849\u00f74='),\n    @('641\u00f72=', '363\u00f78='),\n    @('284\u00f74=', '332\u00f76='),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
